# CA 4.0 files test
# Applies the "natural gas" technology split + new CCS / SMR / hydrogen
# technology rows to the FSCaFoCC sheet, matching the authored diff.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$data  = $wb.Worksheets.Item("Data")
$ws    = $wb.Worksheets.Item("FSCaFoCC")

# ---------------------------------------------------------------------
# 1. Header row: split the old single "Fraction of construction cost
#    (dimensionless)" label in B1 into a two-column header: A1 holds the
#    (italic) unit note, B1 holds the plain description.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Unit: dimensionless"
$ws.Range("A1").Font.Italic = $true

$ws.Range("B1").Value = "Fraction of construction cost"
$ws.Range("B1").Font.Italic = $false

# ---------------------------------------------------------------------
# 2. "natural gas nonpeaker" (row 3) becomes "natural gas steam turbine",
#    and a new "natural gas combined cycle" row is inserted right after
#    it, carrying over the same formula/value.
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "natural gas steam turbine"

$ws.Rows("4:4").Insert()
$ws.Range("A4").Value = "natural gas combined cycle"
$ws.Range("A4").Font.Bold = $true
$ws.Range("B4").Formula = "=Data!B9"

# ---------------------------------------------------------------------
# 3. Append the new technology rows at the bottom of the table (after
#    "municipal solid waste", which is now on row 18 following the
#    insert above). Rows 19-23 reuse the existing bold label style;
#    rows 24-25 (hydrogen techs) use a distinct style (explicit black
#    font colour + vertically centred) matching the authored file.
# ---------------------------------------------------------------------
$newTechs = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor"
)

$row = 19
foreach ($tech in $newTechs) {
    $ws.Cells.Item($row, 1).Value = $tech
    $ws.Cells.Item($row, 1).Font.Bold = $true
    $ws.Cells.Item($row, 2).Value = 0
    $row = $row + 1
}

$hydrogenTechs = @(
    "hydrogen combustion turbine",
    "hydrogen combined cycle"
)

foreach ($tech in $hydrogenTechs) {
    $ws.Cells.Item($row, 1).Value = $tech
    $ws.Cells.Item($row, 1).Font.Bold = $false
    $ws.Cells.Item($row, 1).Font.Color = 0
    $ws.Cells.Item($row, 1).VerticalAlignment = -4108
    $ws.Cells.Item($row, 2).Value = 0
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4. Cosmetic workbook-level changes: colour the FSCaFoCC tab and make
#    it the active/selected sheet (previously "About" was selected).
# ---------------------------------------------------------------------
$ws.Tab.Color = 6437154   # RGB #223962 ~= theme Accent5, darker 50%

$ws.Activate()
$ws.Range("E21").Select()
